$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: blank out the filled-in team name (was shared string "dragon") ---
$ws.Range("B1").Clear() | Out-Null

# --- Row 2: blank out the filled-in group member usernames ---
$ws.Range("B2:C2").Clear() | Out-Null

# --- Row 3: entirely blank out (more filled-in member usernames) so the row disappears ---
$ws.Rows.Item(3).Clear() | Out-Null
$ws.Rows.Item(3).EntireRow.AutoFit() | Out-Null

# --- Row 7: the discussion note goes back to the generic "eg." example text ---
$ws.Range("E7").Value = "eg.
Set up group expectations
Upload and populate the group diary
Fix merge conflict 
Decide on choice of methods for first part"

# --- Row 8: was a fully filled-in meeting entry -> becomes a single italic placeholder ---
$ws.Rows.Item(8).Clear() | Out-Null
$ws.Rows.Item(8).EntireRow.AutoFit() | Out-Null
$ws.Range("A8").Value = "Fill in as needed "
$ws.Range("A8").Font.Italic = $true

# --- Rows 9-11: more filled-in meeting entries, removed entirely ---
$ws.Rows.Item(9).Clear() | Out-Null
$ws.Rows.Item(9).EntireRow.AutoFit() | Out-Null
$ws.Rows.Item(10).Clear() | Out-Null
$ws.Rows.Item(10).EntireRow.AutoFit() | Out-Null
$ws.Rows.Item(11).Clear() | Out-Null
$ws.Rows.Item(11).EntireRow.AutoFit() | Out-Null

# --- Selection back near the top of the sheet on column B ---
$ws.Range("B12").Select() | Out-Null

Write-Output "edits applied"
